$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.072.84"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "1.812.62"

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4626"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3759"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.58%  "

$ws.Range("E9").Value = "  -0.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8635"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("D12").Value = "1.814.40"
$ws.Range("E12").Value = "  +0.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.648"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.392"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07081"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.62%  "

$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008734"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("D21").Value = "27.077.24"
$ws.Range("E21").Value = "  +0.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.335"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.39%  "

$ws.Range("D24").Value = "2.041.34"
$ws.Range("E24").Value = "  +1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.912"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("E28").Value = "  +0.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.267"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08932"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7735"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.18%  "

$ws.Range("E33").Value = "  +0.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.519"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.900"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.32%  "

$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.127"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01957"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05241"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.929"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.235"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5289"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.340"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1676"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.626"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5027"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.672"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06337"
$ws.Range("D51").Style = "Normal"
